# Initial Data File Updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: add two new cells at the end (S11, T11) ---
$ws.Range("S11").Value = "Quincena"
$ws.Range("T11").Value = "Pago Deudas Madre"

# --- Row 13: fill in transaction data (entered before row 12's text) ---
$ws.Range("A13").Value = 43527
$ws.Range("B13").Value = 200
$ws.Range("C13").Value = "Buffete de Carnes"
$ws.Range("D13").Value = "Comida"
$ws.Range("E13").Value = "Gasto"
$ws.Range("F13").Value = "Tarjeta Banamex"
$ws.Range("G13").Value = "Alianza Do Brazil"
$ws.Range("K13").Formula = "=K12-B13"
$ws.Range("L13").Value = 150
$ws.Range("M13").Value = 90
$ws.Range("N13").Formula = "=SUM(K13:M13)"
$ws.Range("O13").Formula = "=N13-4000"

# --- Row 12: fill in transaction data ---
$ws.Range("A12").Value = 43527
$ws.Range("B12").Value = 225
$ws.Range("C12").Value = "Contratación Plan AT&T"
$ws.Range("D12").Value = "Servicios"
$ws.Range("E12").Value = "Gasto"
$ws.Range("F12").Value = "Tarjeta Banamex"
$ws.Range("G12").Value = "AT&T"
$ws.Range("K12").Formula = "=K11-B12"
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = 90
$ws.Range("N12").Formula = "=SUM(K12:M12)"
$ws.Range("O12").Formula = "=N12-4000"
$ws.Range("S12").Value = 5826
$ws.Range("T12").Value = 1500
$ws.Range("U12").Formula = "=S12-T12"

# --- Row 14: fill in transaction data (S14:V14 already existed) ---
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value = 43527
$ws.Range("B14").Value = 20
$ws.Range("C14").Value = "Propina Buffete de Carnes"
$ws.Range("D14").Value = "Comida"
$ws.Range("E14").Value = "Gasto"
$ws.Range("F14").Value = "Efectivo"
$ws.Range("G14").Value = "Alianza Do Brazil"
$ws.Range("K14").Value = 7831.82
$ws.Range("L14").Value = 150
$ws.Range("M14").Formula = "=M13-B14"
$ws.Range("N14").Formula = "=SUM(K14:M14)"
$ws.Range("O14").Formula = "=N14-4000"

# --- Row 15: O15/P15 moved up to S11/T11; update V15 formula ---
$ws.Range("O15").Clear()
$ws.Range("P15").Clear()
$ws.Range("V15").Formula = "=U15+U12"

# --- Row 16: O16/P16/Q16 moved up to S12/T12/U12 ---
$ws.Range("O16").Clear()
$ws.Range("P16").Clear()
$ws.Range("Q16").Clear()

# --- Selection cursor position ---
$ws.Range("Q14").Select()
